# CIERRE 24 NOV 2021
# Fill in the pending credit-collection entries on the NOVIEMBRE 2021 sheet
# (dates received, payer name, amount billed, payment date, amount paid).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("REMISIONES NOVIEMBRE   2021   ")

# Row 5 - DAVID HERRADURA remision already recorded (A5/D5/E5 set);
# register the payment date and amount collected.
$ws.Range("F5").Value = 44511
$ws.Range("G5").Value = 7963

# Row 6 - OBRADOR
$ws.Range("A6").Value = 44510
$ws.Range("D6").Value = "OBRADOR"
$ws.Range("E6").Value = 272
$ws.Range("F6").Value = 44515
$ws.Range("G6").Value = 272

# Row 7 - DAVID HERRADURA
$ws.Range("A7").Value = 44511
$ws.Range("D7").Value = "DAVID HERRADURA"
$ws.Range("E7").Value = 11005
$ws.Range("F7").Value = 44512
$ws.Range("G7").Value = 11005

# Row 8 - DAVID HERRADURA
$ws.Range("A8").Value = 44512
$ws.Range("D8").Value = "DAVID HERRADURA"
$ws.Range("E8").Value = 8335
$ws.Range("F8").Value = 44515
$ws.Range("G8").Value = 8335

# Row 9 - DAVID HERRADURA
$ws.Range("A9").Value = 44515
$ws.Range("D9").Value = "DAVID HERRADURA"
$ws.Range("E9").Value = 4289
$ws.Range("F9").Value = 44516
$ws.Range("G9").Value = 4289

# Row 10 - GUSTAVO
$ws.Range("A10").Value = 44516
$ws.Range("D10").Value = "GUSTAVO"
$ws.Range("E10").Value = 2790
$ws.Range("F10").Value = 44521
$ws.Range("G10").Value = 2790

# Row 11 - OBRADOR
$ws.Range("A11").Value = 44517
$ws.Range("D11").Value = "OBRADOR"
$ws.Range("E11").Value = 196
$ws.Range("F11").Value = 44518
$ws.Range("G11").Value = 196

# Row 12 - GUSTAVO
$ws.Range("A12").Value = 44518
$ws.Range("D12").Value = "GUSTAVO"
$ws.Range("E12").Value = 32652
$ws.Range("F12").Value = 44521
$ws.Range("G12").Value = 32652

# Row 13 - OBRADOR (still unpaid; F/G left blank)
$ws.Range("A13").Value = 44520
$ws.Range("D13").Value = "OBRADOR"
$ws.Range("E13").Value = 152

# Row 14 - MAURO
$ws.Range("A14").Value = 44521
$ws.Range("D14").Value = "MAURO"
$ws.Range("E14").Value = 4649
$ws.Range("F14").Value = 44522
$ws.Range("G14").Value = 4649

# Minor row-height tweak picked up while reviewing row 23
$ws.Rows.Item(23).RowHeight = 15

# Leave the cursor parked where the user left off
$ws.Range("G15").Select()
